{"js": "// Update the date header and the 25 multiplication-table answers to the\n// new day's generated values.\nconst replacements = [\n  [\"2025-08-03 Sunday\", \"2025-08-04 Monday\"],\n  [\"26\u00d715=390\", \"76\u00d765=4940\"],\n  [\"91\u00d727=2457\", \"34\u00d738=1292\"],\n  [\"39\u00d729=1131\", \"62\u00d788=5456\"],\n  [\"25\u00d733=825\", \"89\u00d738=3382\"],\n  [\"15\u00d751=765\", \"89\u00d752=4628\"],\n  [\"32\u00d797=3104\", \"56\u00d749=2744\"],\n  [\"48\u00d782=3936\", \"66\u00d719=1254\"],\n  [\"18\u00d774=1332\", \"22\u00d782=1804\"],\n  [\"31\u00d758=1798\", \"92\u00d726=2392\"],\n  [\"79\u00d726=2054\", \"35\u00d743=1505\"],\n  [\"74\u00d749=3626\", \"97\u00d744=4268\"],\n  [\"74\u00d734=2516\", \"36\u00d743=1548\"],\n  [\"23\u00d731=713\", \"67\u00d745=3015\"],\n  [\"69\u00d734=2346\", \"39\u00d746=1794\"],\n  [\"40\u00d756=2240\", \"15\u00d768=1020\"],\n  [\"72\u00d763=4536\", \"72\u00d794=6768\"],\n  [\"15\u00d797=1455\", \"46\u00d766=3036\"],\n  [\"97\u00d724=2328\", \"75\u00d742=3150\"],\n  [\"72\u00d723=1656\", \"29\u00d729=841\"],\n  [\"13\u00d791=1183\", \"11\u00d718=198\"],\n  [\"86\u00d730=2580\", \"51\u00d738=1938\"],\n  [\"92\u00d776=6992\", \"68\u00d714=952\"],\n  [\"55\u00d750=2750\", \"75\u00d772=5400\"],\n  [\"25\u00d771=1775\", \"50\u00d785=4250\"],\n  [\"95\u00d749=4655\", \"81\u00d740=3240\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: \"${oldText}\"`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date header and the 25 multiplication-table answers to the\n# new day's generated values.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-08-03 Sunday\", \"2025-08-04 Monday\"),\n    @(\"26\u00d715=390\", \"76\u00d765=4940\"),\n    @(\"91\u00d727=2457\", \"34\u00d738=1292\"),\n    @(\"39\u00d729=1131\", \"62\u00d788=5456\"),\n    @(\"25\u00d733=825\", \"89\u00d738=3382\"),\n    @(\"15\u00d751=765\", \"89\u00d752=4628\"),\n    @(\"32\u00d797=3104\", \"56\u00d749=2744\"),\n    @(\"48\u00d782=3936\", \"66\u00d719=1254\"),\n    @(\"18\u00d774=1332\", \"22\u00d782=1804\"),\n    @(\"31\u00d758=1798\", \"92\u00d726=2392\"),\n    @(\"79\u00d726=2054\", \"35\u00d743=1505\"),\n    @(\"74\u00d749=3626\", \"97\u00d744=4268\"),\n    @(\"74\u00d734=2516\", \"36\u00d743=1548\"),\n    @(\"23\u00d731=713\", \"67\u00d745=3015\"),\n    @(\"69\u00d734=2346\", \"39\u00d746=1794\"),\n    @(\"40\u00d756=2240\", \"15\u00d768=1020\"),\n    @(\"72\u00d763=4536\", \"72\u00d794=6768\"),\n    @(\"15\u00d797=1455\", \"46\u00d766=3036\"),\n    @(\"97\u00d724=2328\", \"75\u00d742=3150\"),\n    @(\"72\u00d723=1656\", \"29\u00d729=841\"),\n    @(\"13\u00d791=1183\", \"11\u00d718=198\"),\n    @(\"86\u00d730=2580\", \"51\u00d738=1938\"),\n    @(\"92\u00d776=6992\", \"68\u00d714=952\"),\n    @(\"55\u00d750=2750\", \"75\u00d772=5400\"),\n    @(\"25\u00d771=1775\", \"50\u00d785=4250\"),\n    @(\"95\u00d749=4655\", \"81\u00d740=3240\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $found = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        Write-Output \"NOT FOUND: $oldText\"\n    }\n}\n"}
